# "added basic parser for run management record file"
#
# The sheet "ParserWriter" gets a brand-new row inserted at row 22 that
# documents the new "run management record" (rmr) file parser/writer, which
# pushes every row from the old row 23 onward down by one. Excel re-points
# every cross-sheet formula and the active-sheet/selection bookkeeping
# automatically when a real row is inserted, so we drive the edit with a
# genuine row insert rather than rewriting cells by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ParserWriter")

# --- 1. Insert the new row and fill in the new parser/writer entry -------
$ws.Rows.Item(22).Insert()

$ws.Range("A22").Value = "run management record "
$ws.Range("B22").Value = "rmr"
$ws.Range("C22").Value = 1

# --- 2. Re-point the active sheet / selection -----------------------------
# Previously "PST file" (tab index 2) was active with C31 selected on
# ParserWriter; now ParserWriter itself (tab index 1) is active, selection
# sits on E21, and "PST file" goes back to its plain (unselected) view.
$ws.Activate()
$ws.Range("E21").Select()

# --- 3. Repair the iconSet conditional formatting that covered the now-
#        shifted rows so it keeps tracking the right cells -----------------
# Rows 27:30 became 28:31 and the old blank row 22/23 pair became 22 (new
# data row, not part of the formatting) and 24. Rebuild every rule that
# referenced the old addresses against the new ones, area by area (the
# host's FormatConditions API only reliably accepts single-area ranges).

function Add-IconSetRule {
    param($areas)
    foreach ($addr in $areas) {
        $ws.Range($addr).FormatConditions.AddIconSetCondition() | Out-Null
    }
}

$narrowAreasOld = @("C27:D29", "C30")
$wideAreasOld = @("C27:D29", "C11:D13", "C7:D9", "C22:D23", "D6", "D14:D21", "C30")

$narrowAreasNew = @("C28:D30", "C31")
$wideAreasNew = @("C28:D30", "C11:D13", "C7:D9", "C24:D24", "D6", "D14:D21", "C31", "C22:D22")

# The five rules that applied over the old ranges (in file order): a single
# rule over the "narrow" range, then three more over the "wide" range (the
# last of which was really two stacked cfRules).
$unionOld = $ws.Range($narrowAreasOld[0])
foreach ($a in $narrowAreasOld[1..($narrowAreasOld.Length - 1)]) {
    $unionOld = $excel.Union($unionOld, $ws.Range($a))
}
$staleRules = $unionOld.FormatConditions
$staleCount = $staleRules.Count
for ($i = $staleCount; $i -ge 1; $i--) {
    $staleRules.Item($i).Delete()
}

# priority 37 equivalent
Add-IconSetRule $narrowAreasNew
# priority 137 equivalent
Add-IconSetRule $wideAreasNew
# priority 143 equivalent
Add-IconSetRule $wideAreasNew
# priority 149 + 150 equivalent (two stacked rules over the same cells)
Add-IconSetRule $wideAreasNew
Add-IconSetRule $wideAreasNew

Write-Output "done"
